$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns in order: A, B, D, E, F, G, H, Q, R
# A,B,E,Q,R are numeric; D,F,G,H are strings
$rowsData = @(
    @(112178529, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760450, 7210211),
    @(112178538, 98891, "LC", 222771, "Svart trolldruva", "Actaea spicata", "L.", 760363, 7210127),
    @(112178516, 89331, "LC", 3215, "Rödgul trumpetsvamp", "Craterellus lutescens", "(Fr.) Fr.", 760126, 7210471),
    @(112178514, 102192, "LC", 222412, "Tibast", "Daphne mezereum", "L.", 760068, 7210453),
    @(112178519, 85400, "LC", 1988, "Kryddspindling", "Cortinarius percomis", "Fr.", 760104, 7210466),
    @(112178524, 90826, "LC", 4366, "Skarp dropptaggsvamp", "Hydnellum peckii", "Banker", 760203, 7210420),
    @(112178521, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760097, 7210441),
    @(112178537, 96735, "VU", 220787, "Knärot", "Goodyera repens", "(L.) R. Br.", 760382, 7210147),
    @(112178522, 85448, "NT", 3739, "Persiljespindling", "Cortinarius sulfurinus", "Quél.", 760108, 7210439),
    @(112178539, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760354, 7210135),
    @(112178528, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760519, 7210363),
    @(112178520, 96768, "LC", 219874, "Nattviol", "Platanthera bifolia", "(L.) Rich.", 760092, 7210449),
    @(112178535, 102192, "LC", 222412, "Tibast", "Daphne mezereum", "L.", 760389, 7210155),
    @(112178530, 96735, "VU", 220787, "Knärot", "Goodyera repens", "(L.) R. Br.", 760431, 7210191),
    @(112178518, 85448, "NT", 3739, "Persiljespindling", "Cortinarius sulfurinus", "Quél.", 760120, 7210456),
    @(112178517, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760128, 7210459),
    @(112178532, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760411, 7210179),
    @(112178526, 90814, "LC", 4364, "Dropptaggsvamp", "Hydnellum ferrugineum", "(Fr.:Fr.) P. Karst.", 760256, 7210384),
    @(112178540, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760340, 7210120),
    @(112178515, 90480, "LC", 4769, "Svavelriska", "Lactarius scrobiculatus", "(Scop.:Fr.) Fr.", 760089, 7210467)
)

$startRow = 11
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $row = $rowsData[$i]
    $ws.Cells.Item($r, 1).Value  = $row[0]   # A - Id
    $ws.Cells.Item($r, 2).Value  = $row[1]   # B - Taxonsorteringsordning
    $ws.Cells.Item($r, 4).Value  = $row[2]   # D - Rodlistade
    $ws.Cells.Item($r, 5).Value  = $row[3]   # E - TaxonId
    $ws.Cells.Item($r, 6).Value  = $row[4]   # F - Artnamn
    $ws.Cells.Item($r, 7).Value  = $row[5]   # G - Vetenskapligt namn
    $ws.Cells.Item($r, 8).Value  = $row[6]   # H - Auktor
    $ws.Cells.Item($r, 17).Value = $row[7]   # Q - Ost
    $ws.Cells.Item($r, 18).Value = $row[8]   # R - Nord
}
